# Update the dataText worksheet so the generated sentences "make sense".
# This mirrors a content-only edit: a handful of LatestPeriod / timeTitle
# values get a " data" suffix, and the whole "Forecasted employment ..."
# block (row 21) is replaced with the new "Projected employment growth ..."
# wording, including a brand-new caveat list in column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (LatestPeriod) : append " data" to the existing period text
$ws.Range("B2").Value  = "Oct-Sep 2022 data"
$ws.Range("B3").Value  = "Oct-Sep 2022 data"
$ws.Range("B4").Value  = "Oct-Sep 2022 data"
$ws.Range("B5").Value  = "Oct-Sep 2022 data"
$ws.Range("B6").Value  = "Oct-Sep 2022 data"
$ws.Range("B7").Value  = "Oct-Sep 2022 data"
$ws.Range("B8").Value  = "Oct-Sep 2022 data"
$ws.Range("B9").Value  = "Oct-Sep 2022 data"
$ws.Range("B10").Value = "Dec 2022 data"
$ws.Range("B11").Value = "Mar 2022 data"
$ws.Range("B12").Value = "Dec 2020 - Dec 2021 data"
$ws.Range("B13").Value = "Dec 2020 - Dec 2021 data"
$ws.Range("B14").Value = "AY21/22 data"
$ws.Range("B15").Value = "AY21/22 data"
$ws.Range("B16").Value = "AY21/22 data"
$ws.Range("B17").Value = "AY21/22 data"
$ws.Range("B18").Value = "2021 data"
$ws.Range("B19").Value = "AY20/21 data"
$ws.Range("B20").Value = "AY20/21 data"
$ws.Range("B22").Value = "AY21/22 data"
$ws.Range("B23").Value = "AY21/22 data"

# --- Column H (timeTitle) : reword as a "are ... changing" question
$ws.Range("H2").Value  = "are employment rates changing"
$ws.Range("H3").Value  = "are self-employment rates changing"
$ws.Range("H4").Value  = "are unemployment rates changing"
$ws.Range("H5").Value  = "are inactivity rates changing"
$ws.Range("H6").Value  = "are employment volumes changing"
$ws.Range("H7").Value  = "are self-employment volumes changing"
$ws.Range("H8").Value  = "are unemployment volumes changing"
$ws.Range("H9").Value  = "are inactivity volumes changing"
$ws.Range("H10").Value = "are online job adverts changing"
$ws.Range("H11").Value = "are the number of businesses changing"
$ws.Range("H12").Value = "are business birth rates changing"
$ws.Range("H13").Value = "are business death rates changing"
$ws.Range("H14").Value = "are FE achievement volumes changing"
$ws.Range("H15").Value = "are FE participation volumes changing"
$ws.Range("H16").Value = "are FE achievement rates per 100,000 changing"
$ws.Range("H17").Value = "are FE participation rates per 100,000 changing"
$ws.Range("H18").Value = "are the proportions of people qualified at Level 3 or above changing"
$ws.Range("H19").Value = "are key stage 4 sustained positive destination rates changing"
$ws.Range("H20").Value = "are key stage 5 sustained positive destination rates changing"

# --- Row 21 ("wfEmployment" forecast row) : new wording throughout
$ws.Range("B21").Value = "Growth from 2023 to 2035"
$ws.Range("F21").Value = "<ol>`n  <li>The projections presented in this Workbook are calculated from a number of different data sources, using a variety of econometric and statistical techniques. As a result, precise margins of error cannot be assigned to the estimates. For further details, see the Working Futures Technical Report. </li>`n <li>Industries are based on SIC 2007 codes. </li>`n <li>Time series of the breakdowns can be downloaded in the data download section or in the publication. </li>`n <li>Further breakdowns are available in the published data eg gender, full-time/part-time, as well combined breakdowns. Replacement demand is also available. </li>`n</ol>"
$ws.Range("G21").Value = "Projected employment growth"
$ws.Range("H21").Value = "will year on year employment volume growth change"
$ws.Range("I21").Value = "Employment in"
$ws.Range("J21").Value = "Projected employment growth from 2023 to 2035"
$ws.Range("K21").Value = "projected employment volume changes"
$ws.Range("L21").Value = "projected employment volume change"
$ws.Range("M21").Value = ""

# Row 21 grows taller now that F21 holds a multi-line caveat list (matches
# the autofit Excel itself performs once the wrapped text no longer fits
# the previous row height).
$ws.Rows.Item(21).RowHeight = 203.5

# --- Sheet view: match the new selection recorded in the file
$ws.Range("F21").Select()
